{"js": "// Replace the page title / heading text. This exact string appears twice in\n// the document (the Heading1 title and a bold run further down) and both\n// occurrences get the same replacement, so a single search+replace-all loop\n// covers both.\nconst titleResults = context.document.body.search(\"Play Cash Bunny for Free - Review\", { matchCase: true });\ntitleResults.load(\"items\");\nawait context.sync();\nfor (const r of titleResults.items) {\n  r.insertText(\"Play Cash Bunny for Free - Review and Gameplay\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Simple 1:1 text replacements. Using range.insertText(..., Word.InsertLocation.replace)\n// on the exact matched range (rather than rewriting the whole paragraph) keeps\n// sibling runs -- e.g. the leading empty run in each bullet paragraph -- intact.\nasync function replaceOnce(searchText, newText) {\n  const results = context.document.body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// \"What we like\" bullet list\nawait replaceOnce(\"Fun and simple cartoon farm-theme\", \"Cartoon farm-themed slot game\");\nawait replaceOnce(\"Wild and Bonus symbols add excitement\", \"Wild and Bonus symbols for extra features\");\nawait replaceOnce(\"Two different payline options\", \"Different value levels for symbols\");\nawait replaceOnce(\"Double Spin feature boosts chances of winning\", \"Easy and enjoyable gameplay\");\n\n// \"What we don't like\" bullet list\nawait replaceOnce(\"Basic graphics and no background music\", \"Simple graphics with no animation or background music\");\nawait replaceOnce(\"Low betting limit may not appeal to high rollers\", \"Limited betting options\");\n\n// Meta description (italic run)\nawait replaceOnce(\n  \"Read a review of Cash Bunny, a fun cartoon farm-themed slot game with Wild and Bonus symbols. Try it for free and enjoy its Double Spin feature.\",\n  \"Read our review of Cash Bunny and enjoy a free play of this fun slot game.\"\n);\n", "ps1": "$d = $word.ActiveDocument\n$xmlns = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n\n# Helper: find the paragraph whose full text matches $oldText exactly and\n# rewrite its contents with $newInnerXml (the <w:p>...</w:p> markup for the\n# replacement paragraph). Using InsertXML on the paragraph's own Range keeps\n# sibling runs (e.g. the leading empty <w:r/>) that a plain Range.Text / Find\n# Replace would otherwise collapse away.\nfunction Set-ParagraphXml($oldText, $newParagraphXml) {\n    $target = $null\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text -eq ($oldText + \"`r\")) {\n            $target = $p\n            break\n        }\n    }\n    if ($target -ne $null) {\n        $target.Range.InsertXML($newParagraphXml)\n    }\n}\n\n# Title (Heading 1) -- no empty run in this paragraph\nSet-ParagraphXml \"Play Cash Bunny for Free - Review\" \"<w:p $xmlns><w:pPr><w:pStyle w:val=`\"Heading1`\"/></w:pPr><w:r><w:t>Play Cash Bunny for Free - Review and Gameplay</w:t></w:r></w:p>\"\n\n# \"What we like\" bullet list (each paragraph has a leading empty <w:r/>)\nSet-ParagraphXml \"Fun and simple cartoon farm-theme\" \"<w:p $xmlns><w:pPr><w:pStyle w:val=`\"ListBullet`\"/><w:spacing w:line=`\"240`\" w:lineRule=`\"auto`\"/><w:ind w:left=`\"720`\"/></w:pPr><w:r/><w:r><w:t>Cartoon farm-themed slot game</w:t></w:r></w:p>\"\n\nSet-ParagraphXml \"Wild and Bonus symbols add excitement\" \"<w:p $xmlns><w:pPr><w:pStyle w:val=`\"ListBullet`\"/><w:spacing w:line=`\"240`\" w:lineRule=`\"auto`\"/><w:ind w:left=`\"720`\"/></w:pPr><w:r/><w:r><w:t>Wild and Bonus symbols for extra features</w:t></w:r></w:p>\"\n\nSet-ParagraphXml \"Two different payline options\" \"<w:p $xmlns><w:pPr><w:pStyle w:val=`\"ListBullet`\"/><w:spacing w:line=`\"240`\" w:lineRule=`\"auto`\"/><w:ind w:left=`\"720`\"/></w:pPr><w:r/><w:r><w:t>Different value levels for symbols</w:t></w:r></w:p>\"\n\nSet-ParagraphXml \"Double Spin feature boosts chances of winning\" \"<w:p $xmlns><w:pPr><w:pStyle w:val=`\"ListBullet`\"/><w:spacing w:line=`\"240`\" w:lineRule=`\"auto`\"/><w:ind w:left=`\"720`\"/></w:pPr><w:r/><w:r><w:t>Easy and enjoyable gameplay</w:t></w:r></w:p>\"\n\n# \"What we don't like\" bullet list\nSet-ParagraphXml \"Basic graphics and no background music\" \"<w:p $xmlns><w:pPr><w:pStyle w:val=`\"ListBullet`\"/><w:spacing w:line=`\"240`\" w:lineRule=`\"auto`\"/><w:ind w:left=`\"720`\"/></w:pPr><w:r/><w:r><w:t>Simple graphics with no animation or background music</w:t></w:r></w:p>\"\n\nSet-ParagraphXml \"Low betting limit may not appeal to high rollers\" \"<w:p $xmlns><w:pPr><w:pStyle w:val=`\"ListBullet`\"/><w:spacing w:line=`\"240`\" w:lineRule=`\"auto`\"/><w:ind w:left=`\"720`\"/></w:pPr><w:r/><w:r><w:t>Limited betting options</w:t></w:r></w:p>\"\n\n# Bold title repeated near the end\nSet-ParagraphXml \"Play Cash Bunny for Free - Review\" \"<w:p $xmlns><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Cash Bunny for Free - Review and Gameplay</w:t></w:r></w:p>\"\n\n# Italic meta description\nSet-ParagraphXml \"Read a review of Cash Bunny, a fun cartoon farm-themed slot game with Wild and Bonus symbols. Try it for free and enjoy its Double Spin feature.\" \"<w:p $xmlns><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Cash Bunny and enjoy a free play of this fun slot game.</w:t></w:r></w:p>\"\n"}
